$wb = $excel.ActiveWorkbook

# --- 1) Update the Date value on the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(8, 2).Value = "2026-01-22T09:24:45+00:00"

# --- 2) Insert a new mapping row on "Mapping Table 1" for the new
#        FRCDAResultats.statusCode -> FRDiagnosticReportDocument.status entry ---
$wsMap = $wb.Worksheets.Item("Mapping Table 1")

# Grow the table by one row (row 10), copying the formatting of the last
# existing row (row 9) so the new row matches the table's look.
$wsMap.Range("A9:E9").Copy()
$wsMap.Range("A10:E10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Shift the contents of rows 6-9 down into rows 7-10 (bottom-up, so
# earlier rows aren't overwritten before they're read).
for ($r = 9; $r -ge 6; $r--) {
    for ($c = 1; $c -le 5; $c++) {
        $v = $wsMap.Cells.Item($r, $c).Value2
        $wsMap.Cells.Item($r + 1, $c).Value = $v
    }
}

# Populate the freed-up row 6 with the new status mapping.
$wsMap.Cells.Item(6, 1).Value = "FRCDAResultats.statusCode"
$wsMap.Cells.Item(6, 2).Value = $null
$wsMap.Cells.Item(6, 3).Value = "equivalent"
$wsMap.Cells.Item(6, 4).Value = "FRDiagnosticReportDocument.status"
$wsMap.Cells.Item(6, 5).Value = $null
